$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 1043
$ws.Range("I125").Value = 200
$ws.Range("J125").Value = 1886
$ws.Range("K125").Value = 1800
$ws.Range("L125").Value = 16974
$ws.Range("M125").Value = 660
$ws.Range("N125").Value = -21894
$ws.Range("H132").Value = 440511.4
$ws.Range("I132").Value = 474472.7
$ws.Range("J132").Value = 10335.333
$ws.Range("K132").Value = 1423418.1
$ws.Range("L132").Value = 31005.999
$ws.Range("M132").Value = -1420888.1
$ws.Range("N132").Value = -36065.999
$ws.Range("H137").Value = 4457.5356
$ws.Range("I137").Value = 3872.25
$ws.Range("K137").Value = 11616.75
$ws.Range("M137").Value = -9066.75
$ws.Range("H138").Value = 4291.5854
$ws.Range("J138").Value = 5574.5
$ws.Range("L138").Value = 16723.5
$ws.Range("N138").Value = -27003.5
$ws.Range("H141").Value = 772.85
$ws.Range("I141").Value = 788.8421
$ws.Range("K141").Value = 2366.5263
$ws.Range("M141").Value = 2813.4737

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7470.1055
$ws.Range("I61").Value = 4826.1924
$ws.Range("J61").Value = 13198.583
$ws.Range("K61").Value = 4826.1924
$ws.Range("L61").Value = 13198.583
$ws.Range("M61").Value = -4614.1924
$ws.Range("N61").Value = -13622.583
$ws.Range("H110").Value = 1858.8928
$ws.Range("I110").Value = 1901.3334
$ws.Range("J110").Value = 713
$ws.Range("K110").Value = 1901.3334
$ws.Range("L110").Value = 713
$ws.Range("M110").Value = 143.6666
$ws.Range("N110").Value = -4803
$ws.Range("H132").Value = 340948.47
$ws.Range("I132").Value = 706774.6
$ws.Range("J132").Value = 5607.8057
$ws.Range("K132").Value = 2120323.8
$ws.Range("L132").Value = 16823.4171
$ws.Range("M132").Value = -2117793.8
$ws.Range("N132").Value = -21883.4171
$ws.Range("H136").Value = 7470.1055
$ws.Range("I136").Value = 4826.1924
$ws.Range("J136").Value = 13198.583
$ws.Range("K136").Value = 14478.5772
$ws.Range("L136").Value = 39595.749
$ws.Range("M136").Value = -11928.5772
$ws.Range("N136").Value = -44695.749

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2889.8262
$ws.Range("I99").Value = 1878.6
$ws.Range("K99").Value = 1878.6
$ws.Range("M99").Value = -380.5999999999999
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("N105").ClearContents()
$ws.Range("H134").Value = 498617.2
$ws.Range("I134").Value = 663106.5
$ws.Range("J134").Value = 5149.3335
$ws.Range("K134").Value = 1989319.5
$ws.Range("L134").Value = 15448.0005
$ws.Range("M134").Value = -1986784.5
$ws.Range("N134").Value = -20518.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 318.2857
$ws.Range("J7").Value = 629.6667
$ws.Range("L7").Value = 629.6667
$ws.Range("N7").Value = -855.6667
$ws.Range("H31").Value = 6372.696
$ws.Range("I31").Value = 3605.7273
$ws.Range("K31").Value = 3605.7273
$ws.Range("M31").Value = -3310.7273
$ws.Range("H34").Value = 6372.696
$ws.Range("I34").Value = 3605.7273
$ws.Range("K34").Value = 3605.7273
$ws.Range("M34").Value = -3403.7273
$ws.Range("H58").Value = 444221.22
$ws.Range("I58").Value = 540381.2
$ws.Range("K58").Value = 540381.2
$ws.Range("M58").Value = -540178.2
$ws.Range("H86").Value = 2786.1875
$ws.Range("I86").Value = 2764.7778
$ws.Range("J86").Value = 2813.7144
$ws.Range("K86").Value = 2764.7778
$ws.Range("L86").Value = 2813.7144
$ws.Range("M86").Value = -1641.7778
$ws.Range("N86").Value = -5059.7144
$ws.Range("H89").Value = 2786.1875
$ws.Range("I89").Value = 2764.7778
$ws.Range("J89").Value = 2813.7144
$ws.Range("K89").Value = 13823.889
$ws.Range("L89").Value = 14068.572
$ws.Range("M89").Value = -8207.888999999999
$ws.Range("N89").Value = -25300.572
$ws.Range("H132").Value = 10889897
$ws.Range("I132").Value = 33842.637
$ws.Range("K132").Value = 101527.911
$ws.Range("M132").Value = -98997.91100000001
$ws.Range("H136").Value = 444221.22
$ws.Range("I136").Value = 540381.2
$ws.Range("K136").Value = 1621143.6
$ws.Range("M136").Value = -1618593.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 118121.414
$ws.Range("I14").Value = 118121.414
$ws.Range("K14").Value = 354364.242
$ws.Range("M14").Value = -354191.242
$ws.Range("H23").Value = 296.77777
$ws.Range("I23").Value = 100
$ws.Range("J23").Value = 321.375
$ws.Range("K23").Value = 300
$ws.Range("L23").Value = 964.125
$ws.Range("M23").Value = -65
$ws.Range("N23").Value = -1434.125
$ws.Range("H34").Value = 525
$ws.Range("J34").Value = 250
$ws.Range("L34").Value = 750
$ws.Range("N34").Value = -918
$ws.Range("H132").Value = 3529.5293
$ws.Range("J132").Value = 5339
$ws.Range("L132").Value = 48051
$ws.Range("N132").Value = -53111

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 70071
$ws.Range("J62").Value = 70071
$ws.Range("L62").Value = 70071
$ws.Range("N62").Value = -71443
$ws.Range("H65").Value = 70071
$ws.Range("J65").Value = 70071
$ws.Range("L65").Value = 210213
$ws.Range("N65").Value = -217077
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("N73").ClearContents()
$ws.Range("H97").Value = 1080.375
$ws.Range("I97").Value = 578.56525
$ws.Range("K97").Value = 578.56525
$ws.Range("M97").Value = -82.56524999999999
$ws.Range("H99").Value = 6277.75
$ws.Range("I99").Value = 962
$ws.Range("K99").Value = 962
$ws.Range("M99").Value = 1284
$ws.Range("H132").Value = 2746.0908
$ws.Range("I132").Value = 2349.52
$ws.Range("K132").Value = 7048.559999999999
$ws.Range("M132").Value = -4518.559999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 231948.14
$ws.Range("J20").Value = 255042.45
$ws.Range("L20").Value = 255042.45
$ws.Range("N20").Value = -255494.45
$ws.Range("H35").Value = 1471.0714
$ws.Range("I35").Value = 1471.0714
$ws.Range("K35").Value = 1471.0714
$ws.Range("M35").Value = -1135.0714
$ws.Range("H100").Value = 7946.0625
$ws.Range("I100").Value = 1933.1538
$ws.Range("K100").Value = 1933.1538
$ws.Range("M100").Value = -1392.1538
$ws.Range("H136").Value = 2856.5715
$ws.Range("I136").Value = 2999.4
$ws.Range("K136").Value = 8998.200000000001
$ws.Range("M136").Value = -6448.200000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1603.0667
$ws.Range("I100").Value = 1729.6818
$ws.Range("J100").Value = 1254.875
$ws.Range("K100").Value = 3459.3636
$ws.Range("L100").Value = 2509.75
$ws.Range("M100").Value = -2918.3636
$ws.Range("N100").Value = -3591.75
$ws.Range("H132").Value = 35573050
$ws.Range("I132").Value = 2883270.2
$ws.Range("K132").Value = 8649810.600000001
$ws.Range("M132").Value = -8647280.600000001
$ws.Range("H136").Value = 10857654
$ws.Range("J136").Value = 3281.1667
$ws.Range("L136").Value = 9843.500100000001
$ws.Range("N136").Value = -14943.5001
